$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "56.11") but must remain
# stored as text (matching the source data feed which stores all Price/Volume
# cells as text) - force text format before assigning so Excel does not coerce
# them into numeric cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = '37.206.17'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '2.025.81'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '228.88'
$ws.Range("E5").Value = '  +1.63%  '
$ws.Range("E6").Value = '  +0.87%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '56.11'
$ws.Range("E8").Value = '  +1.78%  '
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").Value = '0.0783'
$ws.Range("E10").Value = '  -0.67%  '
$ws.Range("E11").Value = '  -2.24%  '
$ws.Range("D12").Value = '2.324.84'
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '14.31'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").Value = '20.10'
$ws.Range("E14").Value = '  -1.63%  '
$ws.Range("D15").Value = '5.21'
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '2.025.20'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '37.178.30'
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").Value = '6.18'
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("D20").Value = '69.06'
$ws.Range("E20").Value = '  +0.57%  '
$ws.Range("D21").Value = '0.0₃0819'
$ws.Range("E21").Value = '  -1.47%  '
$ws.Range("D22").Value = '223.40'
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  +1.97%  '
$ws.Range("D25").Value = '2.22'
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("D26").Value = '163.48'
$ws.Range("E26").Value = '  -2.09%  '
$ws.Range("D27").Value = '9.06'
$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("D28").Value = '0.129'
$ws.Range("E28").Value = '  +2.50%  '
$ws.Range("D29").Value = '18.76'
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").Value = '1.32'
$ws.Range("E30").Value = '  -0.85%  '
$ws.Range("D31").Value = '0.117'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("E33").Value = '  -0.84%  '
$ws.Range("D34").Value = '4.47'
$ws.Range("E34").Value = '  +0.46%  '
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("D36").Value = '1.89'
$ws.Range("E36").Value = '  +3.50%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("D39").Value = '5.50'
$ws.Range("E39").Value = '  +2.36%  '
$ws.Range("D40").Value = '1.474.19'
$ws.Range("E40").Value = '  -1.67%  '
$ws.Range("E41").Value = '  -1.75%  '
$ws.Range("D42").Value = '94.74'
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("D43").Value = '2.79'
$ws.Range("E43").Value = '  -2.00%  '
$ws.Range("D44").Value = '0.0915'
$ws.Range("E44").Value = '  -1.54%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '16.32'
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").Value = '4.15'
$ws.Range("E46").Value = '  +15.71%  '
$ws.Range("E47").Value = '  -1.29%  '
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("D49").Value = '7.15'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").Value = '2.211.18'
$ws.Range("E51").Value = '  +0.06%  '
